# Update header text in the first (and only) table.
$d = $word.ActiveDocument

# 1) Replace the header cell texts.
$d.Content.Find.Execute("Case Complete Analysis", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CKD Stages 3a-3b-4", 2)

$d.Content.Find.Execute("Multiple Imputation Analysis", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CKD Stages 3b-4", 2)

# 2) Resize the table grid / columns (widths are in twips -> convert to points: 1 pt = 20 twips).
$table = $d.Tables.Item(1)
$table.Columns.Item(2).Width = 2501 / 20.0
$table.Columns.Item(3).Width = 2183 / 20.0
